$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '31.246.37'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +2.51%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.980.74'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +5.09%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9931'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.74%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.8111'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +71.94%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '252.27'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.40%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9952'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.53%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3423'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +18.58%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.64'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +15.31%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06918'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +6.86%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8467'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +16.79%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08099'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +4.29%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '101.91'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +6.45%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.982.80'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +5.13%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.503'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +6.19%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '275.42'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.18%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '31.246.81'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.52%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.97'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +7.03%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007863'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +5.36%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.236.59'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +4.68%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.678'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +7.78%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9963'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.42%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9949'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.56%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.830'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +8.56%  '

$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1626'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +68.49%  '

$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.634'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +6.24%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.47'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.63%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.57'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.95%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.234'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +18.30%  '

$ws.Range('E30').Value = '  +6.19%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.350'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.07%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.535'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +6.22%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.342'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +4.91%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05177'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +6.64%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.216'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +8.09%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7412'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +7.06%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.772'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.16%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01985'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +5.36%  '

$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.901'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.67%  '

$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.601'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +6.38%  '

$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '78.43'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +4.94%  '

$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4662'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +9.32%  '

$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.076'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +5.79%  '

$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '106.09'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +4.96%  '

$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8545'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.37%  '

$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9961'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.40%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.02'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.91%  '

$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.504'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +7.83%  '

$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '36.39'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.39%  '

$ws.Range('B50').Value = 'Decentraland'
$ws.Range('C50').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4261'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +8.32%  '

$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '934.59'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +3.35%  '
